# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions" - only literal text content
# of column D (Price) and column E (Volume(1h)) changes; everything else
# (labels, links, styles) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.039.38'
$ws.Range('E2').Value = '  +6.55%  '
$ws.Range('D3').Value = '3.012.76'
$ws.Range('E3').Value = '  +3.90%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '584.89'
$ws.Range('E5').Value = '  +2.85%  '
$ws.Range('E6').Value = '  +13.32%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.009.73'
$ws.Range('E8').Value = '  +3.85%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').Value = '  +3.26%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.98'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.159'
$ws.Range('E11').Value = '  +8.58%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  +6.35%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000253'
$ws.Range('E13').Value = '  +9.64%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.97'
$ws.Range('E14').Value = '  +8.15%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '66.027.08'
$ws.Range('D17').Value = '3.511.97'
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.98'
$ws.Range('E18').Value = '  +6.89%  '
$ws.Range('D19').Value = '3.010.99'
$ws.Range('E19').Value = '  +4.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '458.94'
$ws.Range('E20').Value = '  +6.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.96'
$ws.Range('E21').Value = '  +7.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.687'
$ws.Range('E22').Value = '  +5.16%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.39'
$ws.Range('E23').Value = '  +7.48%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.49'
$ws.Range('E24').Value = '  +4.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.32'
$ws.Range('E25').Value = '  +14.35%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.40'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.70'
$ws.Range('E27').Value = '  +6.09%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.12'
$ws.Range('E29').Value = '  +16.07%  '
$ws.Range('E30').Value = '  +15.71%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.62'
$ws.Range('E31').Value = '  +4.69%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0000104'
$ws.Range('E32').Value = '  -5.86%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.00'
$ws.Range('E33').Value = '  +5.52%  '
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +4.43%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.81'
$ws.Range('E37').Value = '  +7.80%  '
$ws.Range('E38').Value = '  +11.77%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.06'
$ws.Range('E39').Value = '  +7.47%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '49.96'
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('E41').Value = '  +14.22%  '
$ws.Range('E42').Value = '  +6.39%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '43.87'
$ws.Range('E43').Value = '  +9.15%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.50'
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '387.16'
$ws.Range('E45').Value = '  +11.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0357'
$ws.Range('E46').Value = '  +6.78%  '
$ws.Range('D47').Value = '2.798.73'
$ws.Range('E47').Value = '  +3.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '135.29'
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.03'
$ws.Range('E51').Value = '  +4.05%  '
